# no-op
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
